# Työaikaraportti - add two new work-log entries (16.2.2024 and 17.2.2024)
# and update the totals ("Yht") row accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The current last row (20) holds the "Yht" (total) row. Insert two new
# rows above it so it becomes row 22, making room for the two new entries
# at rows 20 and 21. (Using row Insert — rather than copying the totals
# row to its new location — keeps its SUM formula/cell correctly wired so
# the later recalculation picks up every contributing row.)
$ws.Rows.Item(20).Insert()
$ws.Rows.Item(20).Insert()

# Seed the formatting of the two new rows from the row above (row 19),
# which carries the same cell styles (date / number / wrapped-text).
$ws.Range("B19:D19").Copy()
$ws.Range("B20:D21").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 20: 16.2.2024, 5 hours
$ws.Range("B20").Value = 45338
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = "Tyylitelyä lisää alan olemaan kunolla tyytyävinen tyyleihin en täysin implementaatioon saatan perehtyä pystyykö paremapaan."

# Row 21: 17.2.2024, 3 hours
$ws.Range("B21").Value = 45339
$ws.Range("C21").Value = 3
$ws.Range("D21").Value = "Säädin tyylejä lisää sain ne aika lailla hyviksi."

# Match the row heights used by the rest of the sheet for these entries.
$ws.Rows.Item(20).RowHeight = 56.25
$ws.Rows.Item(21).RowHeight = 18.75

# Extend the SUM formula in the (now-shifted) totals row to include the
# two new rows.
$ws.Range("C22").Formula = "=SUM(C6:C21)"

# Update the view so the new rows are visible / selected, mirroring the
# author's final cursor position.
$ws.Range("A17").Select()
$ws.Range("D22").Select()
